# Weekly data refresh: shift existing rows down by 3 and insert 3 new
# rows of the latest week's data at the top of the data block (row 16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 16-55 down to 19-58, carrying formatting (incl. date style on
# column D) from the row above, exactly like an Excel UI "Insert Rows".
$ws.Rows("16:18").Insert()

$newDate = Get-Date -Year 2023 -Month 12 -Day 21 -Hour 0 -Minute 0 -Second 0

# New row 16
$ws.Range("A16").Value = 2
$ws.Range("B16").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = $newDate
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = "Fruta"
$ws.Range("G16").Value = 100103
$ws.Range("H16").Value = "Frutos de hueso (carozo)"
$ws.Range("I16").Value = 100103003
$ws.Range("J16").Value = "Damasco"
$ws.Range("K16").Value = "Castle Brite"
$ws.Range("L16").Value = "Especial"
$ws.Range("M16").Value = 60
$ws.Range("N16").Value = 25000
$ws.Range("O16").Value = 26000
$ws.Range("P16").Value = 25500
$ws.Range("Q16").Value = "$/caja 18 kilos"
$ws.Range("R16").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S16").Value = 1417
$ws.Range("T16").Value = 18

# New row 17
$ws.Range("A17").Value = 2
$ws.Range("B17").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C17").Value = "Coquimbo"
$ws.Range("D17").Value = $newDate
$ws.Range("E17").Value = 4
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100103
$ws.Range("H17").Value = "Frutos de hueso (carozo)"
$ws.Range("I17").Value = 100103003
$ws.Range("J17").Value = "Damasco"
$ws.Range("K17").Value = "Castle Brite"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 40
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 23500
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S17").Value = 1306
$ws.Range("T17").Value = 18

# New row 18
$ws.Range("A18").Value = 2
$ws.Range("B18").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C18").Value = "Coquimbo"
$ws.Range("D18").Value = $newDate
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = "Fruta"
$ws.Range("G18").Value = 100103
$ws.Range("H18").Value = "Frutos de hueso (carozo)"
$ws.Range("I18").Value = 100103003
$ws.Range("J18").Value = "Damasco"
$ws.Range("K18").Value = "Castle Brite"
$ws.Range("L18").Value = "Segunda"
$ws.Range("M18").Value = 60
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 21000
$ws.Range("P18").Value = 20500
$ws.Range("Q18").Value = "$/caja 18 kilos"
$ws.Range("R18").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S18").Value = 1139
$ws.Range("T18").Value = 18
